$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new product-code / customer combo (JENNYS re-using gingoso2@gmail.com's email)
$ws.Range("A3").Value = "JENNYS"
$ws.Range("B3").Value = "gingoso2@gmail.com"

# Give B3 its own hyperlink (mailto) just like B1/B2, then restore the plain
# "Hyperlink" cell style (Hyperlinks.Add tends to stamp its own style xf).
$ws.Hyperlinks.Add($ws.Cells.Item(3, 2), "mailto:gingoso2@gmail.com")
$ws.Range("B3").Style = "Hyperlink"

# Move/save the active selection where the user last clicked
[void]$ws.Range("G5").Select()
